$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N:N").Insert()

# The inserted column visually inherits the width of its left neighbour (M)
# in real Excel; reproduce that as closely as the host allows.
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and set its selection
$ws.Activate() | Out-Null
$ws.Range("S7").Select() | Out-Null
